{"js": "// Office.js (Word JavaScript API) edit script.\n// Body is `async (context) => { ... }`.\n//\n// The document currently ends with two empty paragraphs (both \"Normal\"\n// style, lang nb-NO) right before the section break. We turn the first of\n// those into a new \"Heading 1\" section titled \"Video redigeringsprogrammer\"\n// (matching the existing heading convention used throughout this notes\n// document), keep the second empty paragraph as the usual blank spacer\n// line under a heading, and then append two more paragraphs listing the\n// two video editors mentioned in the commit: \"Openshot.org \" and\n// \"Shotcut.org\".\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// The two trailing empty paragraphs just before the end of the body.\nconst count = paragraphs.items.length;\nconst headingPara = paragraphs.items[count - 2];\nconst spacerPara = paragraphs.items[count - 1];\n\n// Turn the first trailing empty paragraph into the new heading, keeping\n// its paragraph-mark language (nb-NO) intact.\nheadingPara.style = \"Heading 1\";\nconst headingParaRange = headingPara.getRange();\nheadingParaRange.languageId = \"nb-NO\";\nconst headingTextRange = headingPara.insertText(\n  \"Video redigeringsprogrammer\",\n  Word.InsertLocation.replace\n);\nheadingTextRange.languageId = \"nb-NO\";\n\n// Leave the second trailing empty paragraph as-is (blank spacer line\n// below the heading), then append the two new body paragraphs after it.\nconst openshotPara = spacerPara.insertParagraph(\"Openshot.org \", Word.InsertLocation.after);\nopenshotPara.insertParagraph(\"Shotcut.org\", Word.InsertLocation.after);\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is the open document.\n#\n# The document currently ends with two empty paragraphs (both \"Normal\"\n# style, lang nb-NO) right before the section break. We turn the first of\n# those into a new \"Heading 1\" section titled \"Video redigeringsprogrammer\"\n# (matching the existing heading convention used throughout this notes\n# document), keep the second empty paragraph as the usual blank spacer\n# line under a heading, and then append two more paragraphs listing the\n# two video editors mentioned in the commit: \"Openshot.org \" and\n# \"Shotcut.org\".\n\n$d = $word.ActiveDocument\n\n$count = $d.Paragraphs.Count\n$headingPara = $d.Paragraphs($count - 1)\n$spacerPara = $d.Paragraphs($count)\n\n# Turn the first trailing empty paragraph into the new heading. Set the\n# text before the style/language so the paragraph-mark's own run\n# properties (the nb-NO language) survive the style switch.\n$headingPara.Range.Text = \"Video redigeringsprogrammer\"\n$headingPara.Style = \"Heading 1\"\n$headingPara.Range.LanguageID = \"nb-NO\"\n\n# Leave the second trailing empty paragraph as-is (blank spacer line below\n# the heading), then append the two new body paragraphs after it.\n$spacerPara.Range.InsertParagraphAfter()\n$openshotPara = $d.Paragraphs($d.Paragraphs.Count)\n$openshotPara.Range.Text = \"Openshot.org \"\n\n$openshotPara.Range.InsertParagraphAfter()\n$shotcutPara = $d.Paragraphs($d.Paragraphs.Count)\n$shotcutPara.Range.Text = \"Shotcut.org\"\n"}
